# Update "想去人数" (want-to-go count) values that changed upstream.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): row 2 and row 3, column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 954
$wsExpo.Range("F3").Value = 1855

# Sheet "全部类型" (All types): row 4 and row 5, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 954
$wsAll.Range("F5").Value = 1855
